# TokenIteratorFieldRewriterSplit: the field-token rewriter now emits a
# dedicated run for each token fragment instead of folding the leading
# "{" (and, for "endlet", the trailing "}") into the neighbouring text
# run. Reproduce that run split here.
#
# Technique: re-assigning a Range's FormattedText to itself is a
# no-visible-effect operation that nonetheless makes Word carve that
# exact character range out into its own run (splitting the run at
# both the start and the end of the range), without adding any new
# run-level formatting (no stray <w:rPr/>).

$d = $word.ActiveDocument

function Isolate-Range($rng) {
    $rng.FormattedText = $rng.FormattedText
}

# --- Change 1 -----------------------------------------------------------
# "{m" (single run, inside the "{m:v}" field) -> "{" + "m" (two runs)
$rng1 = $d.Content.Duplicate
$rng1.Start = 0
$rng1.End = $d.Content.End
if ($rng1.Find.Execute("{m:v}", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0)) {
    $brace = $d.Range($rng1.Start, $rng1.Start + 1)   # the leading "{"
    Isolate-Range $brace
}

# --- Change 2 -------------------------------------------------------------
# "{m:" + "endlet}" (two runs) -> "{" + "m:" + "endlet" + "}" (four runs)
$rng2 = $d.Content.Duplicate
$rng2.Start = 0
$rng2.End = $d.Content.End
if ($rng2.Find.Execute("{m:endlet}", $true, $false, $false, $false, `
                        $false, $true, 1, $false, "", 0)) {
    $start = $rng2.Start
    $end = $rng2.End

    $open = $d.Range($start, $start + 1)              # leading "{"
    Isolate-Range $open

    $close = $d.Range($end - 1, $end)                 # trailing "}"
    Isolate-Range $close
}
